$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 104, pushing existing rows 104..140 down to 105..141.
$ws.Rows.Item(104).Insert()

# Populate the newly inserted row 104 with the new data record.
$ws.Cells.Item(104, 1).Value = 8
$ws.Cells.Item(104, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(104, 3).Value = "Coquimbo"
$ws.Cells.Item(104, 4).Value = 44726
$ws.Cells.Item(104, 5).Value = 4
$ws.Cells.Item(104, 6).Value = 100112001
$ws.Cells.Item(104, 7).Value = "Berenjena"
$ws.Cells.Item(104, 8).Value = "Sin especificar"
$ws.Cells.Item(104, 9).Value = "Primera"
$ws.Cells.Item(104, 10).Value = 520
$ws.Cells.Item(104, 11).Value = 8000
$ws.Cells.Item(104, 12).Value = 9000
$ws.Cells.Item(104, 13).Value = 8500
$ws.Cells.Item(104, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(104, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(104, 16).Value = 170
$ws.Cells.Item(104, 17).Value = 50
$ws.Cells.Item(104, 18).Value = "Hortaliza"
